# Auto-generated edit script applying scheduled market-data refresh
# to the Sargatanas_Profits workbook (one worksheet per crafting job).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 63731924
$ws.Range("I70").Value = 22735306
$ws.Range("J70").Value = 138892380
$ws.Range("K70").Value = 68205918
$ws.Range("L70").Value = 416677140
$ws.Range("M70").Value = -68205648
$ws.Range("N70").Value = -416677680
$ws.Range("H73").Value = 63731924
$ws.Range("I73").Value = 22735306
$ws.Range("J73").Value = 138892380
$ws.Range("K73").Value = 68205918
$ws.Range("L73").Value = 416677140
$ws.Range("M73").Value = -68204982
$ws.Range("N73").Value = -416679012
$ws.Range("H107").Value = 20537116
$ws.Range("I107").Value = 12501127
$ws.Range("J107").Value = 35001896
$ws.Range("K107").Value = 12501127
$ws.Range("L107").Value = 35001896
$ws.Range("M107").Value = -12499207
$ws.Range("N107").Value = -35005736
$ws.Range("H112").Value = 5255.1777
$ws.Range("J112").Value = 5454.442
$ws.Range("L112").Value = 16363.326
$ws.Range("N112").Value = -18579.326
$ws.Range("H132").Value = 1216.3062
$ws.Range("I132").Value = 1168.8889
$ws.Range("J132").Value = 1749.75
$ws.Range("K132").Value = 3506.6667
$ws.Range("L132").Value = 5249.25
$ws.Range("M132").Value = -976.6666999999998
$ws.Range("N132").Value = -10309.25
$ws.Range("H135").Value = 385115.47
$ws.Range("I135").Value = 385115.47
$ws.Range("K135").Value = 3466039.23
$ws.Range("M135").Value = -3463504.23

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4985.5713
$ws.Range("I16").Value = 1099.75
$ws.Range("K16").Value = 1099.75
$ws.Range("M16").Value = -812.75
$ws.Range("H32").Value = 3676.37
$ws.Range("I32").Value = 3785.875
$ws.Range("K32").Value = 3785.875
$ws.Range("M32").Value = -3498.875
$ws.Range("H45").Value = 3474.25
$ws.Range("I45").Value = 3474.25
$ws.Range("K45").Value = 3474.25
$ws.Range("M45").Value = -3097.25
$ws.Range("H61").Value = 5679.5884
$ws.Range("I61").Value = 3029.359
$ws.Range("J61").Value = 14292.833
$ws.Range("K61").Value = 3029.359
$ws.Range("L61").Value = 14292.833
$ws.Range("M61").Value = -2817.359
$ws.Range("N61").Value = -14716.833
$ws.Range("H74").Value = 28182.174
$ws.Range("I74").Value = 38442.25
$ws.Range("J74").Value = 4730.5713
$ws.Range("K74").Value = 38442.25
$ws.Range("L74").Value = 4730.5713
$ws.Range("M74").Value = -37568.25
$ws.Range("N74").Value = -6478.5713
$ws.Range("H77").Value = 28182.174
$ws.Range("I77").Value = 38442.25
$ws.Range("J77").Value = 4730.5713
$ws.Range("K77").Value = 192211.25
$ws.Range("L77").Value = 23652.8565
$ws.Range("M77").Value = -187843.25
$ws.Range("N77").Value = -32388.8565
$ws.Range("H97").Value = 20833772
$ws.Range("I97").Value = 450
$ws.Range("K97").Value = 450
$ws.Range("M97").Value = 46
$ws.Range("H102").Value = 859.9545000000001
$ws.Range("I102").Value = 835.05554
$ws.Range("K102").Value = 835.05554
$ws.Range("M102").Value = 786.94446
$ws.Range("H132").Value = 8260.666999999999
$ws.Range("I132").Value = 8832.583000000001
$ws.Range("J132").Value = 7933.857
$ws.Range("K132").Value = 26497.749
$ws.Range("L132").Value = 23801.571
$ws.Range("M132").Value = -23967.749
$ws.Range("N132").Value = -28861.571
$ws.Range("H136").Value = 5679.5884
$ws.Range("I136").Value = 3029.359
$ws.Range("J136").Value = 14292.833
$ws.Range("K136").Value = 9088.076999999999
$ws.Range("L136").Value = 42878.499
$ws.Range("M136").Value = -6538.076999999999
$ws.Range("N136").Value = -47978.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3275.7273
$ws.Range("I94").Value = 699.3077
$ws.Range("K94").Value = 699.3077
$ws.Range("M94").Value = -248.3077
$ws.Range("H107").Value = 93754210
$ws.Range("I107").Value = 160715150
$ws.Range("K107").Value = 160715150
$ws.Range("M107").Value = -160713230
$ws.Range("H134").Value = 7652.5586
$ws.Range("I134").Value = 3643
$ws.Range("J134").Value = 8692.074000000001
$ws.Range("K134").Value = 10929
$ws.Range("L134").Value = 26076.222
$ws.Range("M134").Value = -8394
$ws.Range("N134").Value = -31146.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10450.833
$ws.Range("I99").Value = 12679.2
$ws.Range("K99").Value = 12679.2
$ws.Range("M99").Value = -11181.2
$ws.Range("H105").Value = 4205515.5
$ws.Range("I105").Value = 5495636.5
$ws.Range("K105").Value = 5495636.5
$ws.Range("M105").Value = -5493889.5
$ws.Range("H126").Value = 10450.833
$ws.Range("I126").Value = 12679.2
$ws.Range("K126").Value = 38037.60000000001
$ws.Range("M126").Value = -35567.60000000001
$ws.Range("H134").Value = 5219.778
$ws.Range("I134").Value = 2280.8386
$ws.Range("J134").Value = 9180.956
$ws.Range("K134").Value = 6842.5158
$ws.Range("L134").Value = 27542.868
$ws.Range("M134").Value = -4307.5158
$ws.Range("N134").Value = -32612.868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1483640.6
$ws.Range("I5").Value = 2667260
$ws.Range("J5").Value = 4116.5
$ws.Range("K5").Value = 8001780
$ws.Range("L5").Value = 12349.5
$ws.Range("M5").Value = -8001668
$ws.Range("N5").Value = -12573.5
$ws.Range("H7").Value = 38.75
$ws.Range("I7").Value = 38.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 116.25
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -4.25
$ws.Range("N7").ClearContents()
$ws.Range("H12").Value = 2941493.5
$ws.Range("J12").Value = 4166790.2
$ws.Range("L12").Value = 12500370.6
$ws.Range("N12").Value = -12500716.6
$ws.Range("H17").Value = 726.5
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 859.8
$ws.Range("K17").Value = 180
$ws.Range("L17").Value = 2579.4
$ws.Range("M17").Value = -11
$ws.Range("N17").Value = -2917.4
$ws.Range("H26").Value = 186.81818
$ws.Range("I26").Value = 52.5
$ws.Range("K26").Value = 157.5
$ws.Range("M26").Value = 130.5
$ws.Range("H32").Value = 126875070
$ws.Range("J32").Value = 145000060
$ws.Range("L32").Value = 435000180
$ws.Range("N32").Value = -435000746
$ws.Range("H34").Value = 5969.0527
$ws.Range("J34").Value = 6640.4116
$ws.Range("L34").Value = 19921.2348
$ws.Range("N34").Value = -20089.2348
$ws.Range("H56").Value = 7749.5
$ws.Range("I56").Value = 7749.5
$ws.Range("K56").Value = 7749.5
$ws.Range("M56").Value = -7219.5
$ws.Range("H124").Value = 3923
$ws.Range("I124").Value = 3923
$ws.Range("K124").Value = 11769
$ws.Range("M124").Value = -6859
$ws.Range("H125").Value = 3999
$ws.Range("I125").Value = 3999
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 11997
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -7077
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 1355.2858
$ws.Range("I126").Value = 1092.3334
$ws.Range("K126").Value = 3277.0002
$ws.Range("M126").Value = 1662.9998
$ws.Range("H129").Value = 1690.4615
$ws.Range("J129").Value = 1879.5555
$ws.Range("L129").Value = 5638.666499999999
$ws.Range("N129").Value = -15638.6665
$ws.Range("H131").Value = 1590
$ws.Range("J131").Value = 1753
$ws.Range("L131").Value = 5259
$ws.Range("N131").Value = -15339
$ws.Range("H135").Value = 1483640.6
$ws.Range("I135").Value = 2667260
$ws.Range("J135").Value = 4116.5
$ws.Range("K135").Value = 24005340
$ws.Range("L135").Value = 37048.5
$ws.Range("M135").Value = -24002805
$ws.Range("N135").Value = -42118.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3949.5
$ws.Range("I102").Value = 3931.5
$ws.Range("J102").Value = 4003.5
$ws.Range("K102").Value = 3931.5
$ws.Range("L102").Value = 4003.5
$ws.Range("M102").Value = -2309.5
$ws.Range("N102").Value = -7247.5
$ws.Range("H132").Value = 8234.388999999999
$ws.Range("I132").Value = 2687.5
$ws.Range("J132").Value = 15168
$ws.Range("K132").Value = 8062.5
$ws.Range("L132").Value = 45504
$ws.Range("M132").Value = -5532.5
$ws.Range("N132").Value = -50564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 38461880
$ws.Range("J55").Value = 474.05884
$ws.Range("L55").Value = 474.05884
$ws.Range("N55").Value = -820.0588399999999
$ws.Range("H122").Value = 4961.85
$ws.Range("I122").Value = 3618.4348
$ws.Range("J122").Value = 6779.4116
$ws.Range("K122").Value = 10855.3044
$ws.Range("L122").Value = 20338.2348
$ws.Range("M122").Value = -8405.304400000001
$ws.Range("N122").Value = -25238.2348
$ws.Range("H132").Value = 10006383
$ws.Range("I132").Value = 20836230
$ws.Range("K132").Value = 62508690
$ws.Range("M132").Value = -62506160
$ws.Range("H136").Value = 12536
$ws.Range("I136").Value = 2800
$ws.Range("K136").Value = 8400
$ws.Range("M136").Value = -5850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 45000
$ws.Range("I15").Value = 45000
$ws.Range("K15").Value = 45000
$ws.Range("M15").Value = -44712
$ws.Range("H136").Value = 38467504
$ws.Range("I136").Value = 111112160
$ws.Range("J136").Value = 8569.058999999999
$ws.Range("K136").Value = 333336480
$ws.Range("L136").Value = 25707.177
$ws.Range("M136").Value = -333333930

